$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.764948844909668
$ws.Range("B1").Value = 1.925398588180542
$ws.Range("C1").Value = 2.191864728927612
$ws.Range("D1").Value = 3.495029449462891
$ws.Range("E1").Value = 2.358052492141724
